# Updated symbol list on Mon Jan 30 09:11:55 UTC 2023 with GitHub Actions
#
# The source sheet stores every data cell (prices, percentages, hour
# counters, coin names, links) as literal TEXT, even though many of the
# values look numeric (e.g. "309.91", "-0.39%", "9"). Excel's COM layer
# auto-coerces a numeric-looking string typed into `.Value` into a real
# number, which would change both the stored type and the cell style. To
# keep the text type (and the original "no explicit style" look) we:
#   1. assign the value with a leading apostrophe, which forces Excel to
#      store it as text (adds a transient "quote prefix" style), then
#   2. immediately copy the style from an untouched plain-text cell (a
#      "Coin" cell in the same row) back onto the cell, which clears the
#      quote-prefix style without touching the stored text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainTextValue($cell, $value, $styleTemplate) {
    $cell.Value = "'" + $value
    $cell.Style = $styleTemplate.Style
}

# Per-row updates taken from the diff. Columns omitted for a row are
# unchanged (e.g. rows 27-38 keep "--" / "--%" in D/E, so only G moves).
$rowUpdates = @(
    @{ Row = 2; D = '309.91'; E = '-0.39%'; G = '9' },
    @{ Row = 3; D = '38.27'; E = '-2.54%'; G = '9' },
    @{ Row = 4; D = '5.139'; E = '0.21%'; G = '9' },
    @{ Row = 5; D = '0.07998'; E = '-1.43%'; G = '9' },
    @{ Row = 6; B = 'FTXToken'; C = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; D = '2.068'; E = '2.86%'; G = '9' },
    @{ Row = 7; B = 'GateToken'; C = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; D = '4.478'; E = '5.56%'; G = '9' },
    @{ Row = 8; D = '8.308'; E = '2.04%'; G = '9' },
    @{ Row = 9; D = '3.114'; E = '-7.06%'; G = '9' },
    @{ Row = 10; D = '0.9409'; E = '1.46%'; G = '9' },
    @{ Row = 11; D = '0.1305'; E = '-8.14%'; G = '9' },
    @{ Row = 12; D = '0.1922'; E = '-0.52%'; G = '9' },
    @{ Row = 13; D = '0.08966'; E = '-0.94%'; G = '9' },
    @{ Row = 14; D = '0.03477'; E = '-1.70%'; G = '9' },
    @{ Row = 15; D = '0.09719'; E = '-1.04%'; G = '9' },
    @{ Row = 16; D = '0.001415'; E = '1.03%'; G = '9' },
    @{ Row = 17; D = '0.006757'; E = '11.83%'; G = '9' },
    @{ Row = 18; D = '3.584'; E = '-5.14%'; G = '9' },
    @{ Row = 19; D = '0.3465'; E = '0.34%'; G = '9' },
    @{ Row = 20; D = '0.1294'; E = '-1.43%'; G = '9' },
    @{ Row = 21; D = '5.042'; E = '8.12%'; G = '9' },
    @{ Row = 22; D = '0.2532'; E = '4.26%'; G = '9' },
    @{ Row = 23; D = '0.04376'; E = '-0.08%'; G = '9' },
    @{ Row = 24; D = '0.001249'; E = '1.57%'; G = '9' },
    @{ Row = 25; D = '0.004683'; E = '-2.29%'; G = '9' },
    @{ Row = 26; D = '0.0003588'; E = '175.66%'; G = '9' },
    @{ Row = 27; G = '9' },
    @{ Row = 28; G = '9' },
    @{ Row = 29; G = '9' },
    @{ Row = 30; G = '9' },
    @{ Row = 31; G = '9' },
    @{ Row = 32; G = '9' },
    @{ Row = 33; G = '9' },
    @{ Row = 34; G = '9' },
    @{ Row = 35; G = '9' },
    @{ Row = 36; G = '9' },
    @{ Row = 37; G = '9' },
    @{ Row = 38; G = '9' },
    @{ Row = 39; D = '0.02180'; E = '1.59%'; G = '9' },
    @{ Row = 40; D = '0.05167'; E = '0.92%'; G = '9' },
    @{ Row = 41; D = '0.007639'; E = '2.50%'; G = '9' },
    @{ Row = 42; D = '0.01001'; E = '2.04%'; G = '9' },
    @{ Row = 43; D = '0.1389'; E = '2.06%'; G = '9' },
    @{ Row = 44; D = '0.002045'; E = '-3.64%'; G = '9' },
    @{ Row = 45; D = '0.009127'; E = '6.08%'; G = '9' },
    @{ Row = 46; D = '0.00006683'; E = '4.63%'; G = '9' },
    @{ Row = 47; D = '0.00000000756'; E = '0.62%'; G = '9' },
    @{ Row = 48; D = '0.003024'; E = '17.76%'; G = '9' },
    @{ Row = 49; D = '0.001210'; E = '20.80%'; G = '9' },
    @{ Row = 50; D = '0.00002117'; E = '0.62%'; G = '9' },
    @{ Row = 51; D = '0.0002016'; E = '0.62%'; G = '9' }
)

foreach ($item in $rowUpdates) {
    $r = $item.Row
    $styleTemplate = $ws.Cells.Item($r, 2)

    if ($item.ContainsKey('B')) {
        $ws.Cells.Item($r, 2).Value = $item.B
    }
    if ($item.ContainsKey('C')) {
        $ws.Cells.Item($r, 3).Value = $item.C
    }
    if ($item.ContainsKey('D')) {
        Set-PlainTextValue $ws.Cells.Item($r, 4) $item.D $styleTemplate
    }
    if ($item.ContainsKey('E')) {
        Set-PlainTextValue $ws.Cells.Item($r, 5) $item.E $styleTemplate
    }
    if ($item.ContainsKey('G')) {
        Set-PlainTextValue $ws.Cells.Item($r, 7) $item.G $styleTemplate
    }
}
